$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update the price list (column D, rows 31-34)
$ws.Range("D31").Value = 1515.938
$ws.Range("D32").Value = 1998.938
$ws.Range("D33").Value = 2529.188
$ws.Range("D34").Value = 2852.063
